$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) values remain stored as text, matching the
# original inline-string cell type, rather than being auto-converted to numbers.
$ws.Range("D2:D50").NumberFormat = "@"

# --- Column D (Price) text updates ---
$ws.Range("D2").Value = '236.18'
$ws.Range("D3").Value = '21.81'
$ws.Range("D4").Value = '5.358'
$ws.Range("D5").Value = '0.05590'
$ws.Range("D6").Value = '6.480'
$ws.Range("D7").Value = '3.354'
$ws.Range("D8").Value = '0.7996'
$ws.Range("D9").Value = '1.042'
$ws.Range("D10").Value = '0.01167'
$ws.Range("D11").Value = '0.1389'
$ws.Range("D12").Value = '0.07308'
$ws.Range("D13").Value = '0.03171'
$ws.Range("D14").Value = '0.02953'
$ws.Range("D15").Value = '0.09243'
$ws.Range("D16").Value = '0.001681'
$ws.Range("D17").Value = '3.254'
$ws.Range("D18").Value = '0.04771'
$ws.Range("D19").Value = '0.006226'
$ws.Range("D20").Value = '0.005054'
$ws.Range("D22").Value = '0.0001502'
$ws.Range("D23").Value = '0.0003723'
$ws.Range("D24").Value = '3.947'
$ws.Range("D40").Value = '0.04120'
$ws.Range("D41").Value = '0.007128'
$ws.Range("D42").Value = '0.003505'
$ws.Range("D43").Value = '0.1037'
$ws.Range("D44").Value = '0.008788'
$ws.Range("D45").Value = '0.00005444'
$ws.Range("D47").Value = '0.6762'
$ws.Range("D48").Value = '0.03516'
$ws.Range("D49").Value = '0.00002103'
$ws.Range("D50").Value = '0.01012'

# --- Other column updates (Coin name, Link, Volume label) ---
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("E10").Value = '9OneONEBestin24h'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("E11").Value = '10WazirXWRX'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'MCDex'
$ws.Range("C17").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("E17").Value = '16MCDexMCB'
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("E20").Value = '19HotbitTokenHTB'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("E43").Value = '42BKEXTokenBKK'
